$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 42613.761006944442
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = 38
$ws.Range("E8").Value = 58
$ws.Range("F8").Value = 25
$ws.Range("G8").Value = 29313
$ws.Range("H8").Value = 17168
$ws.Range("I8").Value = 952
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 134
$ws.Range("L8").Value = 15
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = "Named"

# Row 9
$ws.Range("A9").Value = 42613.890694444446
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = 42
$ws.Range("E9").Value = 55
$ws.Range("F9").Value = 25
$ws.Range("G9").Value = 17125
$ws.Range("H9").Value = 16197
$ws.Range("I9").Value = 917
$ws.Range("J9").Value = 185
$ws.Range("K9").Value = 142
$ws.Range("L9").Value = 15
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = "Named"

# Row 10
$ws.Range("A10").Value = 42614.887523148151
$ws.Range("B10").Value = 36
$ws.Range("C10").Value = 59
$ws.Range("D10").Value = 38
$ws.Range("E10").Value = 59
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 22275
$ws.Range("H10").Value = 12708
$ws.Range("I10").Value = 698
$ws.Range("J10").Value = 161
$ws.Range("K10").Value = 104
$ws.Range("L10").Value = 17
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = "Named"

# Row 11
$ws.Range("A11").Value = 42615.886793981481
$ws.Range("B11").Value = 42
$ws.Range("C11").Value = 62
$ws.Range("D11").Value = 36
$ws.Range("E11").Value = 62
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 19561
$ws.Range("H11").Value = 15641
$ws.Range("I11").Value = 860
$ws.Range("J11").Value = 195
$ws.Range("K11").Value = 115
$ws.Range("L11").Value = 22
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = "Named"
